$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Credentials")

$rows = @(
    @{ Row = 36; A = "firwl@gmail.com"; B = "1Aemst205^"; C = "Karthik";    D = "Agarwal D";   E = "Karthik Agarwal D";   F = "," },
    @{ Row = 37; A = "bmzim@gmail.com"; B = "T4bekl736*"; C = "Bhadraksh";  D = "Kapoor D";     E = "Bhadraksh Kapoor D";  F = "," },
    @{ Row = 38; A = "dclob@gmail.com"; B = "KCn62l343!"; C = "Shubha";    D = "Pothuvaal D";  E = "Shubha Pothuvaal D";  F = "," }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
